$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: fix label text and strip the bold/centered header styling ---
$ws.Range("C1").Value = "DEVISI"
$ws.Rows(1).ClearFormats()
$ws.Rows(1).AutoFit()

# --- Birth-date text corrections ---
$ws.Range("B4").Value = "Bontokape, 01  - Oktober 1997"
$ws.Range("B5").Value = "Jember, 16-Oktober 1998"
$ws.Range("B10").Value = "Denpasar, 04 Oktober 1993"

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection moves to C1 ---
$ws.Range("C1").Select() | Out-Null
